$d = $word.ActiveDocument

# --- Part 1: "Loyer maximum mensuel :" -> Jinja if/else block ---------
# Locate the run containing "Loyer maximum mensuel :" (bold run at the
# start of the "Loyer" paragraph).
$find1 = $d.Content
$find1.Find.ClearFormatting()
$ok1 = $find1.Find.Execute("Loyer maximum mensuel :")
if (-not $ok1) {
    throw "Could not find 'Loyer maximum mensuel :'"
}

$p1Xml = '<?xml version="1.0" standalone="yes"?>' +
'<?mso-application progid="Word.Document"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body><w:p>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">{% if </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>programme.is_foyer</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>programme.is_residence</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:proofErr w:type="gramStart"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>%}Redevance</w:t></w:r>' +
'<w:proofErr w:type="gramEnd"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">{% </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>else</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> %}</w:t></w:r>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Loyer</w:t></w:r>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{% endif %}</w:t></w:r>' +
'<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> maximum mensuel :</w:t></w:r>' +
'</w:p></w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

# Clear the old text, then splice the new run sequence in its place.
$find1.Text = ""
$ins1 = $d.Range($find1.Start, $find1.Start)
$ins1.InsertXML($p1Xml)

# --- Part 2: "{{ loyer" + "_m2" -> single run "{{ loyer_m2" -----------
$findA = $d.Content
$findA.Find.ClearFormatting()
$okA = $findA.Find.Execute("{{ loyer")
if (-not $okA) {
    throw "Could not find '{{ loyer'"
}
$afterA = $findA.End

$findB = $d.Range($afterA, $afterA + 3)
if ($findB.Text -ne "_m2") {
    throw "Unexpected text after '{{ loyer': [$($findB.Text)]"
}

# Merge "{{ loyer" and "_m2" into a single run "{{ loyer_m2" that
# carries the formatting of the original "{{ loyer" run (no more
# gramStart/gramEnd proofErr wrapping around it).
$mergedRange = $d.Range($findA.Start, $findB.End)
$mergedRange.Text = "{{ loyer_m2"
